# Applies the OOXML diff: shrinks the "label" column of the three
# client-data tables by 425 dxa (21.25 pt) and shifts the tables right
# by the same amount so their right edge stays put; also sets bottom
# vertical alignment on the small spacer cell in the header row of the
# first of those tables.

$d = $word.ActiveDocument

# Word's Tables collection is 1-based and document-order; Tables(1) is
# the unrelated auto-width table earlier in the doc, so the three
# affected tables are Tables(2..4).
$tableIndexes = @(2, 3, 4)

foreach ($ti in $tableIndexes) {
    $t = $d.Tables.Item($ti)

    # <w:tblW w:w="9639"/> -> <w:tblW w:w="9214"/>  (481.95pt -> 460.7pt)
    $t.PreferredWidth = 460.7

    # <w:tblInd w:w="-567"/> -> <w:tblInd w:w="-142"/>  (-28.35pt -> -7.1pt)
    $t.Rows.LeftIndent = -7.1
}

# --- Table 2 (doc order) : "DADOS PESSOAIS" table ---
$t1 = $d.Tables.Item(2)

# Header row, gridSpan=4 cell: tcW 7938 -> 7513 (396.9pt -> 375.65pt)
$t1.Rows.Item(1).Cells.Item(1).Width = 375.65

# Header row, trailing spacer cell (tcW 1701, unchanged width) gains
# <w:vAlign w:val="bottom"/>
$t1.Rows.Item(1).Cells.Item(2).VerticalAlignment = 3

# Row 2 "Nome Completo" first cell: tcW 2835 -> 2410 (141.75pt -> 120.5pt)
$t1.Rows.Item(2).Cells.Item(1).Width = 120.5

# Row 3 "Naturalidade" first cell: tcW 2835 -> 2410
$t1.Rows.Item(3).Cells.Item(1).Width = 120.5

# Row 4 "Email" cell (gridSpan=2): tcW 5245 -> 4820 (262.25pt -> 241pt)
$t1.Rows.Item(4).Cells.Item(1).Width = 241

# Row 5 "CEP" first cell: tcW 2835 -> 2410
$t1.Rows.Item(5).Cells.Item(1).Width = 120.5

# --- Table 3 (doc order) : "DADOS FUNCIONAIS" table ---
$t2 = $d.Tables.Item(3)

# Header row, gridSpan=3 cell: tcW 7938 -> 7513
$t2.Rows.Item(1).Cells.Item(1).Width = 375.65

# Row 2 "Fonte Pagadora" first cell: tcW 3404 -> 2979 (170.2pt -> 148.95pt)
$t2.Rows.Item(2).Cells.Item(1).Width = 148.95

# --- Table 4 (doc order) : "DADOS DA OPERACAO" table ---
$t3 = $d.Tables.Item(4)

# Header row, gridSpan=2 cell: tcW 8080 -> 7655 (404pt -> 382.75pt)
$t3.Rows.Item(1).Cells.Item(1).Width = 382.75

# All remaining data rows: first cell tcW 5480 -> 5055 (274pt -> 252.75pt)
for ($r = 2; $r -le $t3.Rows.Count; $r++) {
    $t3.Rows.Item($r).Cells.Item(1).Width = 252.75
}

Write-Host "table width/indent edits applied"
